$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$co = $ws.ChartObjects(1)
$chart = $co.Chart
try { Write-Host $chart.ChartData } catch { Write-Host "no ChartData: $_" }
try { $chart.ChartData.Activate() } catch { Write-Host "no Activate: $_" }
try { Write-Host $chart.PlotBy } catch { Write-Host "no PlotBy: $_" }
